# Updated cryptos list on Tue Jun 20 20:35:11 UTC 2023 with GitHub Actions
#
# Writes new price/volume snapshot values into the existing cryptos table.
# Every Price/Link/Coin cell on this sheet is stored as literal text (the
# OOXML has t="inlineStr"/shared-string cells, never numeric cells), so we
# force text entry the same way a human typist would force Excel to keep a
# number-looking value as text: a leading apostrophe. Afterwards we restore
# the cell's style to "Normal" so no incidental NumberFormat/style id is
# left behind on cells that had none before.

function Set-CellText {
    param($ws, $addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-CellText $ws "D2" "27.972.13"
Set-CellText $ws "E2" "  +4.77%  "

# Row 3 - Ethereum
Set-CellText $ws "D3" "1.780.92"
Set-CellText $ws "E3" "  +3.11%  "

# Row 4 - TetherUSD
Set-CellText $ws "D4" "0.9998"
Set-CellText $ws "E4" "  +0.25%  "

# Row 5 - BNB
Set-CellText $ws "D5" "244.13"
Set-CellText $ws "E5" "  +0.93%  "

# Row 6 - USDC
Set-CellText $ws "E6" "  +0.24%  "

# Row 7 - XRP
Set-CellText $ws "D7" "0.4921"
Set-CellText $ws "E7" "  +0.00%  "

# Row 8 - Cardano
Set-CellText $ws "D8" "0.2672"
Set-CellText $ws "E8" "  +2.13%  "

# Row 9 - Dogecoin
Set-CellText $ws "D9" "0.06255"
Set-CellText $ws "E9" "  +0.35%  "

# Row 10 - WrappedEther
Set-CellText $ws "D10" "1.784.66"
Set-CellText $ws "E10" "  +3.32%  "

# Row 11 - Solana
Set-CellText $ws "D11" "16.36"
Set-CellText $ws "E11" "  +3.13%  "

# Row 12 - TRON
Set-CellText $ws "D12" "0.07043"
Set-CellText $ws "E12" "  +0.68%  "

# Row 13 - Polygon
Set-CellText $ws "D13" "0.6269"
Set-CellText $ws "E13" "  +2.55%  "

# Row 14 - Polkadot
Set-CellText $ws "D14" "4.632"
Set-CellText $ws "E14" "  +2.84%  "

# Row 15 - Litecoin
Set-CellText $ws "D15" "80.09"
Set-CellText $ws "E15" "  +3.77%  "

# Row 16/17 - Dai and WrappedBTC swap position (WrappedBTC moves above Dai)
Set-CellText $ws "B16" "WrappedBTC"
Set-CellText $ws "C16" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-CellText $ws "D16" "27.936.30"
Set-CellText $ws "E16" "  +5.39%  "

Set-CellText $ws "B17" "Dai"
Set-CellText $ws "C17" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText $ws "D17" "0.9998"
Set-CellText $ws "E17" "  +0.20%  "

# Row 18 - BinanceUSD
Set-CellText $ws "D18" "0.9995"
Set-CellText $ws "E18" "  +0.23%  "

# Row 19 - ShibaInu
Set-CellText $ws "D19" "0.000007217"
Set-CellText $ws "E19" "  +0.04%  "

# Row 20 - Avalanche
Set-CellText $ws "D20" "11.97"
Set-CellText $ws "E20" "  +4.80%  "

# Row 21 - Wrapped liquid staked Ether 2.0
Set-CellText $ws "D21" "2.009.98"
Set-CellText $ws "E21" "  +3.26%  "

# Row 22 - Uniswap
Set-CellText $ws "E22" "  +2.98%  "

# Row 23 - Cosmos
Set-CellText $ws "D23" "8.704"

# Row 24 - Chainlink
Set-CellText $ws "D24" "5.229"
Set-CellText $ws "E24" "  +2.42%  "

# Row 25 - Monero
Set-CellText $ws "D25" "141.77"
Set-CellText $ws "E25" "  +2.56%  "

# Row 26 - EthereumClassic
Set-CellText $ws "D26" "15.74"
Set-CellText $ws "E26" "  +2.40%  "

# Row 27 - LidoDAOToken
Set-CellText $ws "D27" "1.861"
Set-CellText $ws "E27" "  +6.23%  "

# Row 28 - BitcoinCash
Set-CellText $ws "D28" "109.55"
Set-CellText $ws "E28" "  +3.16%  "

# Row 29 - Toncoin
Set-CellText $ws "D29" "1.389"
Set-CellText $ws "E29" "  +0.11%  "

# Row 30 - InternetComputer(DFINITY)
Set-CellText $ws "D30" "4.199"
Set-CellText $ws "E30" "  +7.24%  "

# Row 31 - Stellar
Set-CellText $ws "D31" "0.08295"
Set-CellText $ws "E31" "  +3.82%  "

# Row 32 - Filecoin
Set-CellText $ws "D32" "3.795"
Set-CellText $ws "E32" "  +3.40%  "

# Row 33 - Hedera
Set-CellText $ws "D33" "0.04878"
Set-CellText $ws "E33" "  +8.60%  "

# Row 34 - ARBITRUM
Set-CellText $ws "D34" "1.071"
Set-CellText $ws "E34" "  +6.78%  "

# Row 35 - HuobiToken
Set-CellText $ws "E35" "  +0.17%  "

# Row 36 - ImmutableX
Set-CellText $ws "D36" "0.6507"
Set-CellText $ws "E36" "  +3.88%  "

# Row 37 - TrustWalletToken
Set-CellText $ws "D37" "0.9489"
Set-CellText $ws "E37" "  +1.45%  "

# Row 38 - MXToken
Set-CellText $ws "D38" "2.585"
Set-CellText $ws "E38" "  +7.11%  "

# Row 39 - RenderToken
Set-CellText $ws "D39" "2.046"
Set-CellText $ws "E39" "  +1.54%  "

# Row 40 - FraxShare
Set-CellText $ws "D40" "5.978"
Set-CellText $ws "E40" "  +6.88%  "

# Row 41 - VeChain
Set-CellText $ws "D41" "0.01552"
Set-CellText $ws "E41" "  +2.55%  "

# Row 42 - PaxDollar
Set-CellText $ws "D42" "0.9992"
Set-CellText $ws "E42" "  +0.17%  "

# Row 43 - Quant
Set-CellText $ws "D43" "100.01"
Set-CellText $ws "E43" "  +0.48%  "

# Row 44 - TheSandbox
Set-CellText $ws "D44" "0.3986"
Set-CellText $ws "E44" "  +3.20%  "

# Row 45 - Aptos
Set-CellText $ws "D45" "7.186"
Set-CellText $ws "E45" "  +3.91%  "

# Row 46 - Algorand
Set-CellText $ws "D46" "0.1203"
Set-CellText $ws "E46" "  +3.66%  "

# Row 47 - Cronos
Set-CellText $ws "E47" "  +0.62%  "

# Row 48 - EnergySwap
Set-CellText $ws "D48" "8.011"
Set-CellText $ws "E48" "  +2.47%  "

# Row 49 - NEARProtocol
Set-CellText $ws "D49" "1.299"
Set-CellText $ws "E49" "  +5.16%  "

# Row 50 - Elrond
Set-CellText $ws "D50" "30.66"
Set-CellText $ws "E50" "  +1.02%  "

# Row 51 - Aave
Set-CellText $ws "D51" "52.92"
Set-CellText $ws "E51" "  +2.34%  "
